$d = $word.ActiveDocument

# --- Paragraph 1: split "AQUATOX 4.0, differences in scientific formula from Release 3.2"
#     into three runs: "...formula" | "s" | " from Release 3.2"
#     (net effect: "formula from" -> "formulas from", with new run boundaries)
$target1 = "AQUATOX 4.0, differences in scientific formula from Release 3.2"
$r1 = $d.Range(0, $target1.Length)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>AQUATOX 4.0, differences in scientific formula</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> from Release 3.2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Paragraph 2: merge "Internal Nutrients " + "update" runs into a single run
$d.Content.Find.Execute("Internal Nutrients update", $false, $false, $false, $false, $false, $true, 1, $false, "Internal Nutrients update", 2)
